$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New word-level rows (text, x, y, width, height) replacing the old
# phrase-level rows 2-5.
$data = @(
    @("segregation ", 524, 644, 103, 23),
    @("Battle ", 621, 694, 55, 23),
    @("of ", 676, 694, 22, 23),
    @("Normandy. ", 698, 694, 102, 23),
    @("historically ", 477, 719, 99, 23),
    @("black ", 576, 719, 50, 23),
    @("college ", 626, 719, 63, 23),
    @("NAACP's ", 477, 968.1999999999999, 82, 23),
    @("boycotts ", 397, 993.1999999999999, 76, 23),
    @("school ", 572, 1059.6, 59, 23),
    @("integration, ", 631, 1059.6, 104, 23),
    @("civil ", 268, 1151, 40, 23),
    @("rights ", 308, 1151, 54, 23),
    @("activists ", 362, 1151, 74, 23)
)

$rowIndex = 2
foreach ($item in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $item[0]
    $ws.Cells.Item($rowIndex, 2).Value = $item[1]
    $ws.Cells.Item($rowIndex, 3).Value = $item[2]
    $ws.Cells.Item($rowIndex, 4).Value = $item[3]
    $ws.Cells.Item($rowIndex, 5).Value = $item[4]
    $rowIndex++
}
